$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh: update price (D) and volume-change (E) columns,
# plus two coin-row swaps (PEPE/Maker and Stellar/ApeXProtocol) per source diff.

$ws.Range("D2").Value = "69.843.76"
$ws.Range("E2").Value = "  +0.32%  "

$ws.Range("D3").Value = "3.508.78"
$ws.Range("E3").Value = "  -0.55%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.44%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "196.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.51%  "

$ws.Range("E7").Value = "  +1.32%  "

$ws.Range("E8").Value = "  -0.11%  "

$ws.Range("E9").Value = "  -2.17%  "

$ws.Range("E10").Value = "  +2.06%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.11"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.55%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000300"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.28%  "

$ws.Range("E13").Value = "  +1.28%  "

$ws.Range("D14").Value = "4.063.00"
$ws.Range("E14").Value = "  -0.66%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "602.69"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.05%  "

$ws.Range("D16").Value = "70.016.71"
$ws.Range("E16").Value = "  +0.46%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.02"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.09%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.03%  "

$ws.Range("D19").Value = "3.509.12"
$ws.Range("E19").Value = "  +0.22%  "

$ws.Range("E20").Value = "  +0.83%  "

$ws.Range("E21").Value = "  +0.48%  "

$ws.Range("E22").Value = "  +3.53%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "103.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.92%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.08%  "

$ws.Range("E25").Value = "  -2.59%  "

$ws.Range("E26").Value = "  +3.40%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.15%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.19%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.29%  "

$ws.Range("E30").Value = "  +28.32%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.28%  "

$ws.Range("E32").Value = "  +4.55%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.115"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.97%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.36%  "

$ws.Range("B35").Value = "Maker"
$ws.Range("C35").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D35").Value = "3.734.06"
$ws.Range("E35").Value = "  +5.64%  "

$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").Value = "0.0₃0816"
$ws.Range("E36").Value = "  +5.03%  "

$ws.Range("E37").Value = "  -4.42%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.01%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.393"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.56%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.78"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.28%  "

$ws.Range("E41").Value = "  +1.28%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "500.17"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.73%  "

$ws.Range("E43").Value = "  +0.36%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0458"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.07%  "

$ws.Range("E45").Value = "  -3.00%  "

$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.140"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.73%  "

$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.89%  "

$ws.Range("E48").Value = "  +0.33%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.69"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.08%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000244"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.39%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "130.34"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.41%  "

